$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.406.28"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.89"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.24"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3766"
$ws.Range("E7").Value = "  +2.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.94"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3423"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07654"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.158"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.27"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.029"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.953"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.572.63"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001134"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.21"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.85"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.430"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.397.77"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.706"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.31"
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.20"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.046"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.40"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.748.24"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.182"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.007"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9971"
$ws.Range("E34").Value = "  -4.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.01"
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08595"
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2319"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06590"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.327"
$ws.Range("E40").Value = "  +6.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.459"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6439"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.52"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.07"
$ws.Range("E45").Value = "  -3.53%  "
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6001"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.305"
$ws.Range("E48").Value = "  +7.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.088"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.67"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07332"
$ws.Range("E51").Value = "  +0.47%  "
